# Remove file upload functionality
#
# The underlying data-collection workbook logs one row per day for each of
# the four device/channel sheets (DE_LFT_#1, DE_LFT_#2, DE_PLT_#1, DE_PLT_#2).
# The last recorded row (row 89) is duplicated into a new row 90 with an
# updated timestamp (one day later) - everything else on the row stays the
# same as the prior entry.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Duplicate the last data row (row 89) down into the new row 90,
    # carrying over every column's value/format exactly as-is.
    $ws.Range("A89:I89").Copy($ws.Range("A90:I90"))

    # The new entry's timestamp is one day after the previous one.
    $ws.Cells.Item(90, 1).Value = 45876.43542824074
}
